# Swap the presentation's applied theme colour scheme from the
# "Integral / Red Violet" palette over to the "Office" palette
# (ppt/theme/theme2.xml -- the theme actually referenced by the
# presentation and its slide master / all slide layouts).
#
# PowerPoint's ThemeColorScheme exposes the 12 theme colour slots in
# the same fixed order they appear in <a:clrScheme>: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink. .RGB uses the standard VBA RGB()
# packing (R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
